# CloudWalk Technical Case - Anti-Fraud Solution
# 1) Add <w:lang w:val="de-DE"/> to the paragraph-mark rPr and the run rPr of
#    the "Prototype Code Snippets" heading paragraph.
# 2) Same addition for the "@app.post("/score")" paragraph right after it.
# 3) Fix the split word "hold_for_revie" + "w rates." -> "hold_for_review" + " rates."
#    while keeping the two separate runs (bold "hold_for_review" run, then the
#    regular " rates." run) intact.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Get-ParagraphByText($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

# --- Edit 1: "Prototype Code Snippets" paragraph ---
$p1 = Get-ParagraphByText $d "Prototype Code Snippets*"
$xml1 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="463DD577" w14:textId="77777777" w:rsidR="003C5997" w:rsidRPr="00EB20DB" w:rsidRDefault="003C5997" w:rsidP="00DF5670">' + `
    '<w:pPr><w:spacing w:before="120" w:after="120" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="de-DE"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EB20DB"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="de-DE"/></w:rPr><w:t>Prototype Code Snippets</w:t></w:r>' + `
    '</w:p>'
$p1.Range.InsertXML($xml1)

# --- Edit 2: "@app.post(""/score"")" paragraph ---
$p2 = Get-ParagraphByText $d '*@app.post("/score")*'
$xml2 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="7B0ED8C7" w14:textId="77777777" w:rsidR="00EB20DB" w:rsidRPr="00EB20DB" w:rsidRDefault="00EB20DB" w:rsidP="00DF5670">' + `
    '<w:pPr><w:spacing w:before="120" w:after="120" w:line="276" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="de-DE"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EB20DB"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="de-DE"/></w:rPr><w:t>@app.post("/score")</w:t></w:r>' + `
    '</w:p>'
$p2.Range.InsertXML($xml2)

# --- Edit 3: fix "hold_for_revie" / "w rates." split ---
$p3 = Get-ParagraphByText $d "*Alert unusual spikes in hold_for_review rates.*"
$xml3 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="5EEC76B5" w14:textId="44BA568E" w:rsidR="00EB20DB" w:rsidRPr="00EB20DB" w:rsidRDefault="00EB20DB" w:rsidP="00DF5670">' + `
    '<w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="13"/></w:numPr><w:spacing w:before="120" w:after="120" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="00EB20DB"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Alert unusual spikes in </w:t></w:r>' + `
    '<w:r w:rsidRPr="00EB20DB"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>hold_for_review</w:t></w:r>' + `
    '<w:r w:rsidRPr="00EB20DB"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> rates.</w:t></w:r>' + `
    '</w:p>'
$p3.Range.InsertXML($xml3)

Write-Host "Done."
